# Add concept tags to column F for several dataset rows.
# Order mirrors the sequence in which the new shared strings were authored
# (one-way, loglinear;multinomial;zeros, ca, glm, loglinear, glm;logistic,
# loglinear;logit;2x2) so the resulting sharedStrings.xml ordering matches.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F11").Value = "one-way"
$ws.Range("F13").Value = "one-way"
$ws.Range("F38").Value = "one-way"
$ws.Range("F5").Value = "loglinear;multinomial;zeros"
$ws.Range("F7").Value = "ca"
$ws.Range("F12").Value = "glm"
$ws.Range("F14").Value = "loglinear"
$ws.Range("F15").Value = "glm;logistic"
$ws.Range("F2").Value = "loglinear;logit;2x2"

# Update the selected cell to match the author's final cursor position.
$ws.Range("F18").Select()
